# Applies the FFXIV leve-profit recalculation updates described in the commit.
# Each block targets one worksheet/row; values are written cell-by-cell to match
# the refreshed currentAveragePrice / LevePrice / LeveProfit figures exactly.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2247.0833
$ws.Range("I40").Value = 2361
$ws.Range("J40").Value = 2133.1667
$ws.Range("K40").Value = 2361
$ws.Range("L40").Value = 2133.1667
$ws.Range("M40").Value = -2186
$ws.Range("N40").Value = -2483.1667
# Row 76
$ws.Range("H76").Value = 6188.8335
$ws.Range("I76").Value = 4656.8
$ws.Range("J76").Value = 13849
$ws.Range("K76").Value = 4656.8
$ws.Range("L76").Value = 13849
$ws.Range("M76").Value = -4341.8
$ws.Range("N76").Value = -14479
# Row 79
$ws.Range("H79").Value = 6188.8335
$ws.Range("I79").Value = 4656.8
$ws.Range("J79").Value = 13849
$ws.Range("K79").Value = 4656.8
$ws.Range("L79").Value = 13849
$ws.Range("M79").Value = -3564.8
$ws.Range("N79").Value = -16033
# Row 113
$ws.Range("H113").Value = 2368.2727
$ws.Range("I113").Value = 2460.1
$ws.Range("K113").Value = 2460.1
$ws.Range("M113").Value = 793.9000000000001
# Row 116
$ws.Range("H116").Value = 4188.0386
$ws.Range("I116").Value = 4851.9414
$ws.Range("J116").Value = 2934
$ws.Range("K116").Value = 4851.9414
$ws.Range("L116").Value = 2934
$ws.Range("M116").Value = -1409.9414
$ws.Range("N116").Value = -9818
# Row 137
$ws.Range("H137").Value = 7829.1665
$ws.Range("I137").Value = 8359.48
$ws.Range("J137").Value = 5177.6
$ws.Range("K137").Value = 25078.44
$ws.Range("L137").Value = 15532.8
$ws.Range("M137").Value = -22528.44
$ws.Range("N137").Value = -20632.8
# Row 138
$ws.Range("H138").Value = 2538.5754
$ws.Range("I138").Value = 1535.2413
$ws.Range("J138").Value = 3199.8635
$ws.Range("K138").Value = 4605.7239
$ws.Range("L138").Value = 9599.5905
$ws.Range("M138").Value = 534.2761
$ws.Range("N138").Value = -19879.5905

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 13316.667
$ws.Range("J46").Value = 13316.667
$ws.Range("L46").Value = 13316.667
$ws.Range("N46").Value = -13954.667
# Row 74
$ws.Range("H74").Value = 2014.8235
$ws.Range("I74").Value = 2014.8235
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2014.8235
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1140.8235
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 2014.8235
$ws.Range("I77").Value = 2014.8235
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10074.1175
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5706.1175
$ws.Range("N77").ClearContents()
# Row 110
$ws.Range("H110").Value = 1923.439
$ws.Range("I110").Value = 1858.4117
$ws.Range("K110").Value = 1858.4117
$ws.Range("M110").Value = 186.5882999999999
# Row 122
$ws.Range("H122").Value = 60330.332
$ws.Range("J122").Value = 60330.332
$ws.Range("L122").Value = 180990.996
$ws.Range("N122").Value = -185890.996
# Row 125
$ws.Range("H125").Value = 115062.664
$ws.Range("J125").Value = 115062.664
$ws.Range("L125").Value = 115062.664
$ws.Range("N125").Value = -124902.664
# Row 132
$ws.Range("H132").Value = 2227.9697
$ws.Range("I132").Value = 2227.9697
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6683.909100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4153.909100000001
$ws.Range("N132").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2982.9524
$ws.Range("I134").Value = 2586.5789
$ws.Range("J134").Value = 6748.5
$ws.Range("K134").Value = 7759.736699999999
$ws.Range("L134").Value = 20245.5
$ws.Range("M134").Value = -5224.736699999999
$ws.Range("N134").Value = -25315.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2760
$ws.Range("I31").Value = 1898.1818
$ws.Range("K31").Value = 1898.1818
$ws.Range("M31").Value = -1603.1818
# Row 34
$ws.Range("H34").Value = 2760
$ws.Range("I34").Value = 1898.1818
$ws.Range("K34").Value = 1898.1818
$ws.Range("M34").Value = -1696.1818
# Row 132
$ws.Range("H132").Value = 6281.2163
$ws.Range("I132").Value = 2038.1923
$ws.Range("K132").Value = 6114.5769
$ws.Range("M132").Value = -3584.5769
# Row 134
$ws.Range("H134").Value = 2753.15
$ws.Range("I134").Value = 2706.162
$ws.Range("J134").Value = 3332.6667
$ws.Range("K134").Value = 8118.485999999999
$ws.Range("L134").Value = 9998.000100000001
$ws.Range("M134").Value = -5583.485999999999
$ws.Range("N134").Value = -15068.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 931095.1
$ws.Range("J9").Value = 200848.6
$ws.Range("L9").Value = 602545.8
$ws.Range("N9").Value = -602993.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 41678000
$ws.Range("I80").Value = 83350340
$ws.Range("J80").Value = 5666.3335
$ws.Range("K80").Value = 83350340
$ws.Range("L80").Value = 5666.3335
$ws.Range("M80").Value = -83349342
$ws.Range("N80").Value = -7662.3335
# Row 83
$ws.Range("H83").Value = 41678000
$ws.Range("I83").Value = 83350340
$ws.Range("J83").Value = 5666.3335
$ws.Range("K83").Value = 416751700
$ws.Range("L83").Value = 28331.6675
$ws.Range("M83").Value = -416746708
$ws.Range("N83").Value = -38315.6675
# Row 107
$ws.Range("H107").Value = 6441.7646
$ws.Range("I107").Value = 705.7143
$ws.Range("J107").Value = 10457
$ws.Range("K107").Value = 705.7143
$ws.Range("L107").Value = 10457
$ws.Range("M107").Value = 1214.2857
$ws.Range("N107").Value = -14297
# Row 109
$ws.Range("H109").Value = 32285
$ws.Range("J109").Value = 32285
$ws.Range("L109").Value = 32285
$ws.Range("N109").Value = -34365
# Row 132
$ws.Range("H132").Value = 2507.5454
$ws.Range("I132").Value = 2026.1428
$ws.Range("J132").Value = 3350
$ws.Range("K132").Value = 6078.428400000001
$ws.Range("L132").Value = 10050
$ws.Range("M132").Value = -3548.428400000001
$ws.Range("N132").Value = -15110

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 5000569
$ws.Range("I16").Value = 6250448
$ws.Range("K16").Value = 6250448
$ws.Range("M16").Value = -6250278
# Row 122
$ws.Range("H122").Value = 4930
$ws.Range("I122").Value = 11745.75
$ws.Range("J122").Value = 2451.5454
$ws.Range("K122").Value = 35237.25
$ws.Range("L122").Value = 7354.6362
$ws.Range("M122").Value = -32787.25
$ws.Range("N122").Value = -12254.6362
# Row 132
$ws.Range("H132").Value = 6582.846
$ws.Range("I132").Value = 3111.8
$ws.Range("J132").Value = 11316.091
$ws.Range("K132").Value = 9335.400000000001
$ws.Range("L132").Value = 33948.273
$ws.Range("M132").Value = -6805.400000000001
$ws.Range("N132").Value = -39008.273
